# edit.ps1 - apply the daily-progress report edits described by the diff:
#   1) Table cell "22" (Score) -> split into two runs "2" + "0" (same rPr)
#   2) Table cell course-title run merge: "WordPress Training for Beginners " +
#      "From" + " Scratch" (wrapped in proofErr tags) -> single run
#      "WordPress Training for Beginners From Scratch" (proofErr removed)
#   3) Table cell "Eduonix" -> drop the spellStart/spellEnd proofErr wrapper
#   4) Big paragraph: merge the "The course title is ... CMS WordPress. "
#      run-fragments (with gramStart/gramEnd/spellStart/spellEnd proofErr
#      wrappers) into a single run, proofErr removed.
#
# Because the engine always normalises/merges adjacent runs that resolve to
# identical formatting when edited through Range.Text/InsertAfter, the only
# way to reproduce the exact run layout the diff expects (including the
# "22" -> "2"+"0" split) is to replace each affected paragraph's contents
# wholesale via Range.InsertXML with hand-built OOXML that mirrors the
# paragraph's existing rsid/paraId metadata and run formatting.

$d = $word.ActiveDocument

function Replace-ParagraphXML($findText, $innerXml, $paraAttrs) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($findText) | Out-Null
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body><w:p ' + $paraAttrs + '>' + $innerXml + '</w:p></w:body>' +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) Score cell: "22" -> "2" + "0" (two runs, identical Times New Roman rPr)
# ---------------------------------------------------------------------
$rPr1 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$inner1 = '<w:pPr>' + $rPr1 + '</w:pPr>' +
    '<w:r>' + $rPr1 + '<w:t>2</w:t></w:r>' +
    '<w:r>' + $rPr1 + '<w:t>0</w:t></w:r>'
Replace-ParagraphXML "22" $inner1 'w14:paraId="7F081291" w14:textId="07487BB3" w:rsidR="00FD482B" w:rsidRPr="00323384" w:rsidRDefault="00F05952" w:rsidP="00323384"'

# ---------------------------------------------------------------------
# 2) Course-title cell: merge 3 runs (incl. gramStart/gramEnd proofErr) into one
# ---------------------------------------------------------------------
$rPr2 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$inner2 = '<w:pPr>' + $rPr2 + '</w:pPr>' +
    '<w:r><w:t>WordPress Training for Beginners From Scratch</w:t></w:r>'
Replace-ParagraphXML "WordPress Training for Beginners From Scratch" $inner2 'w14:paraId="1D7C0094" w14:textId="48DE2519" w:rsidR="006D2F12" w:rsidRPr="008365AB" w:rsidRDefault="00373E8A" w:rsidP="008365AB"'

# ---------------------------------------------------------------------
# 3) Certificate-provider cell: "Eduonix" -> drop spellStart/spellEnd wrapper
# ---------------------------------------------------------------------
$rPr3 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$inner3 = '<w:pPr>' + $rPr3 + '</w:pPr>' +
    '<w:r>' + $rPr3 + '<w:t>Eduonix</w:t></w:r>'
Replace-ParagraphXML "Eduonix" $inner3 'w14:paraId="6DDD34B1" w14:textId="114A6010" w:rsidR="006D2F12" w:rsidRPr="006D2F12" w:rsidRDefault="00373E8A" w:rsidP="00853208"'

# ---------------------------------------------------------------------
# 4) Body paragraph: merge "The course title is ... CMS WordPress. " run
#    fragments (incl. gramStart/gramEnd + spellStart/spellEnd proofErr) into
#    a single run.
# ---------------------------------------------------------------------
$rPr4 = '<w:rPr><w:rFonts w:ascii="Arial Black" w:hAnsi="Arial Black"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$bodyText = 'The course title is “WordPress Training for Beginners From Scratch” from the Eduonix, one of the best online course content provider. The course will teach some easy techniques to design a website with world famous CMS WordPress. '
$inner4 = '<w:pPr>' + $rPr4 + '</w:pPr>' +
    '<w:r>' + $rPr4 + '<w:t xml:space="preserve">' + $bodyText + '</w:t></w:r>'

$r1 = $d.Content.Duplicate
$r1.Find.Execute("The course title is") | Out-Null
$startPos = $r1.Start
$r2 = $d.Content.Duplicate
$r2.Find.Execute("CMS WordPress. ") | Out-Null
$endPos = $r2.End
$full = $d.Range($startPos, $endPos)
$pkg4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body><w:p w14:paraId="3FE16B5A" w14:textId="02866A4E" w:rsidR="004B1A4F" w:rsidRDefault="004B1A4F" w:rsidP="005F19EF">' + $inner4 + '</w:p></w:body>' +
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($pkg4)

Write-Output "done"
